$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix laporan pendapatan export excel ---
# The bottom "summary" rows (Total Penjualan, Total Diskon Produk, Total Diskon
# Nota, Total Pengeluaran, Total Transfer, Modal Usaha, Laba Bersih) are no
# longer part of this export template: unmerge the label cells, wipe their
# text, and drop the now-unused alignment so the rows fall back to plain
# (left column) / bold-vertical-centered (the two rows that were bold)
# formatting. The amount column (D) also loses its stray top border so the
# whole block renders as a clean, un-annotated form.

$labelRows = 13..19

foreach ($r in $labelRows) {
    $rng = $ws.Range("B$r`:C$r")
    if ($rng.MergeCells) {
        $rng.UnMerge()
    }
    $rng.ClearContents()
    $rng.HorizontalAlignment = 1   # xlGeneral - drop the explicit left/right alignment
}

# D17/D19 previously carried a top border to underline the "Total Transfer"
# and "Laba Bersih" rows; remove it now that the labels are gone.
foreach ($r in 16, 17, 18, 19) {
    $ws.Range("D$r").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> xlLineStyleNone
}

# Leave the cursor where the edit finished, like the source change.
$ws.Range("D18").Select()
